# Adding objective IQ functionality
# Insert a new "ParentIndex" column between the existing "IsComposite" (D)
# and "Filename" (old E, now F) columns, and populate it with 1 for every
# data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old "Filename" column (and everything from it) one column to
# the right, opening up column E for the new data.
$ws.Columns("E:E").Insert()

# Header for the new column.
$ws.Range("E1").Value = "ParentIndex"

# Every row gets a ParentIndex of 1.
$ws.Range("E2:E12").Value = 1

# Match the new column's width to its neighbor (IsComposite, column D).
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# Leave the selection where the author left off.
$ws.Range("E12").Select() | Out-Null
